$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for 9c0d3966-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-27 06:42:33"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 9c0d3966-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-27 06:42:29"
$wsZhCn.Range("K4").Value = "2016-08-27 06:42:47"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 9c0d3966-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-27 06:42:33"
$wsDeDe.Range("K4").Value = "2016-08-27 06:42:53"
